# Added missing Cool_Low_Building to Sets_Fuel and Set_filter_file
#
# 1) Technology_selection: append 14 new "HLDH_*" technologies after the
#    existing last row (148 -> 149..162).
# 2) Fuel_selection: insert a new "Cool_Low_Building" fuel at row 20,
#    pushing the existing rows (and the dependent data-validation ranges)
#    down by one.

$wb = $excel.ActiveWorkbook

# --- Technology_selection: append new technologies -------------------------
$wsTech = $wb.Worksheets.Item("Technology_selection")

$newTechnologies = @(
    "HLDH_Biomass_Boiler",
    "HLDH_Biomass_CHP",
    "HLDH_Coal_Boiler",
    "HLDH_Coal_CHP",
    "HLDH_Oil_Boiler",
    "HLDH_Gas_Boiler",
    "HLDH_Gas_CHP",
    "HLDH_Geothermal",
    "HLDH_Solar_Thermal",
    "HLDH_WasteToEnergy_Boiler",
    "HLDH_WasteToEnergy_CHP",
    "HLDH_Heatpump",
    "HLDH_ExcessHeat",
    "HLDH_Electric_Boiler"
)

$firstNewTechRow = 149
for ($i = 0; $i -lt $newTechnologies.Count; $i++) {
    $row = $firstNewTechRow + $i
    $wsTech.Cells.Item($row, 1).Value = $newTechnologies[$i]
    $wsTech.Cells.Item($row, 2).Value = 1
}

$wsTech.Range("B153").Select() | Out-Null

# --- Fuel_selection: insert Cool_Low_Building at row 20 --------------------
$wsFuel = $wb.Worksheets.Item("Fuel_selection")

$wsFuel.Rows.Item(20).Insert() | Out-Null
$wsFuel.Cells.Item(20, 1).Value = "Cool_Low_Building"
$wsFuel.Cells.Item(20, 2).Value = 1

$wsFuel.Columns.Item(1).ColumnWidth = 17.42

# Page setup for the Technology_selection sheet (portrait / A4-ish "9")
$wsTech.PageSetup.PaperSize = 9
$wsTech.PageSetup.Orientation = 1

# --- Best-effort view/selection bookkeeping ---------------------------------
$wsYear = $wb.Worksheets.Item("Year_selection")
$wsYear.Range("D20").Select() | Out-Null

$wsFuel.Activate() | Out-Null
$wsFuel.Range("F9").Select() | Out-Null

Write-Output "Inserted Cool_Low_Building into Fuel_selection and appended HLDH_* technologies to Technology_selection"
